# "Generate Report for Archive"
#
# The localization status changed from "Ready for handoff" to
# "In Translation" for the 81a5b25e-... source file. Update every cell
# that shows this status (Overview!E2/F2, zh-cn!C2, de-de!C2) and then
# re-fit the now-narrower Status columns, the same way Excel would after
# a user edits a cell and re-sizes the column to the new content.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# Narrow the Status columns to fit the shorter text (was sized for
# "Ready for handoff", now only needs to fit "In Translation").
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
